$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,16
$row2[0,0] = 3
$row2[0,1] = 1
$row2[0,2] = 21.84976866666667
$row2[0,3] = 65.549306
$row2[0,4] = 0.05020018890879543
$row2[0,5] = 0.05020018890879543
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 8.970048
$row2[0,9] = 26.910144
$row2[0,10] = 0.487108783009476
$row2[0,11] = 0.4871087830094759
$row2[0,12] = 195.993473728896
$row2[0,13] = 1763.941263560064
$row2[0,14] = 0.02445295292620913
$row2[0,15] = 0.02445295292620913
$ws.Range("E2:T2").Value = $row2

$row3 = New-Object "object[,]" 1,16
$row3[0,0] = 3
$row3[0,1] = 1
$row3[0,2] = 21.84976866666667
$row3[0,3] = 65.549306
$row3[0,4] = 0.05020018890879543
$row3[0,5] = 0.05020018890879543
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 9.012070666666666
$row3[0,9] = 27.036212
$row3[0,10] = 0.489390778604016
$row3[0,11] = 0.489390778604016
$row3[0,12] = 196.9116592743191
$row3[0,13] = 1772.204933468872
$row3[0,14] = 0.02456750953614408
$row3[0,15] = 0.02456750953614408
$ws.Range("E3:T3").Value = $row3

$row4 = New-Object "object[,]" 1,16
$row4[0,0] = 3
$row4[0,1] = 1
$row4[0,2] = 21.84976866666667
$row4[0,3] = 65.549306
$row4[0,4] = 0.05020018890879543
$row4[0,5] = 0.05020018890879543
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 0.4327576666666667
$row4[0,9] = 1.298273
$row4[0,10] = 0.02350043838650813
$row4[0,11] = 0.02350043838650813
$row4[0,12] = 9.45565490539311
$row4[0,13] = 85.100894148538
$row4[0,14] = 0.001179726446442216
$row4[0,15] = 0.001179726446442216
$ws.Range("E4:T4").Value = $row4

$row5 = New-Object "object[,]" 1,16
$row5[0,0] = 3
$row5[0,1] = 1
$row5[0,2] = 385.0524703333334
$row5[0,3] = 1155.157411
$row5[0,4] = 0.8846641374295412
$row5[0,5] = 0.8846641374295412
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 8.970048
$row5[0,9] = 26.910144
$row5[0,10] = 0.487108783009476
$row5[0,11] = 0.4871087830094759
$row5[0,12] = 3453.939141408576
$row5[0,13] = 31085.45227267719
$row5[0,14] = 0.4309276713554316
$row5[0,15] = 0.4309276713554316
$ws.Range("E5:T5").Value = $row5

$row6 = New-Object "object[,]" 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 385.0524703333334
$row6[0,3] = 1155.157411
$row6[0,4] = 0.8846641374295412
$row6[0,5] = 0.8846641374295412
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 9.012070666666666
$row6[0,9] = 27.036212
$row6[0,10] = 0.489390778604016
$row6[0,11] = 0.489390778604016
$row6[0,12] = 3470.120073018571
$row6[0,13] = 31231.08065716713
$row6[0,14] = 0.4329464710196934
$row6[0,15] = 0.4329464710196934
$ws.Range("E6:T6").Value = $row6

$row7 = New-Object "object[,]" 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 385.0524703333334
$row7[0,3] = 1155.157411
$row7[0,4] = 0.8846641374295412
$row7[0,5] = 0.8846641374295412
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 0.4327576666666667
$row7[0,9] = 1.298273
$row7[0,10] = 0.02350043838650813
$row7[0,11] = 0.02350043838650813
$row7[0,12] = 166.6344086056892
$row7[0,13] = 1499.709677451203
$row7[0,14] = 0.02078999505441629
$row7[0,15] = 0.02078999505441629
$ws.Range("E7:T7").Value = $row7

$row8 = New-Object "object[,]" 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 28.350479
$row8[0,3] = 85.05143699999999
$row8[0,4] = 0.06513567366166337
$row8[0,5] = 0.06513567366166337
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 8.970048
$row8[0,9] = 26.910144
$row8[0,10] = 0.487108783009476
$row8[0,11] = 0.4871087830094759
$row8[0,12] = 254.305157452992
$row8[0,13] = 2288.746417076928
$row8[0,14] = 0.03172815872783522
$row8[0,15] = 0.03172815872783522
$ws.Range("E8:T8").Value = $row8

$row9 = New-Object "object[,]" 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 28.350479
$row9[0,3] = 85.05143699999999
$row9[0,4] = 0.06513567366166337
$row9[0,5] = 0.06513567366166337
$row9[0,6] = 3
$row9[0,7] = 1
$row9[0,8] = 9.012070666666666
$row9[0,9] = 27.036212
$row9[0,10] = 0.489390778604016
$row9[0,11] = 0.489390778604016
$row9[0,12] = 255.4965201818493
$row9[0,13] = 2299.468681636644
$row9[0,14] = 0.03187679804817853
$row9[0,15] = 0.03187679804817853
$ws.Range("E9:T9").Value = $row9

$row10 = New-Object "object[,]" 1,16
$row10[0,0] = 3
$row10[0,1] = 1
$row10[0,2] = 28.350479
$row10[0,3] = 85.05143699999999
$row10[0,4] = 0.06513567366166337
$row10[0,5] = 0.06513567366166337
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 0.4327576666666667
$row10[0,9] = 1.298273
$row10[0,10] = 0.02350043838650813
$row10[0,11] = 0.02350043838650813
$row10[0,12] = 12.26888714092233
$row10[0,13] = 110.419984268301
$row10[0,14] = 0.00153071688564962
$row10[0,15] = 0.00153071688564962
$ws.Range("E10:T10").Value = $row10

Write-Output "applied updates to rows 2-10"